$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 88.281043
$ws.Cells.Item(2, 8).Value = 264.843129
$ws.Cells.Item(2, 9).Value = 0.5104595351890647
$ws.Cells.Item(2, 10).Value = 0.6037079925017727
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 88.281043
$ws.Cells.Item(2, 14).Value = 264.843129
$ws.Cells.Item(2, 15).Value = 0.5104595351890647
$ws.Cells.Item(2, 16).Value = 0.6037079925017727
$ws.Cells.Item(2, 17).Value = 7793.542553167848
$ws.Cells.Item(2, 18).Value = 70141.88297851064
$ws.Cells.Item(2, 19).Value = 0.260568937065436
$ws.Cells.Item(2, 20).Value = 0.3644633402105204

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 88.281043
$ws.Cells.Item(3, 8).Value = 264.843129
$ws.Cells.Item(3, 9).Value = 0.5104595351890647
$ws.Cells.Item(3, 10).Value = 0.6037079925017727
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 3.371552333333333
$ws.Cells.Item(3, 14).Value = 10.114657
$ws.Cells.Item(3, 15).Value = 0.01949502382905625
$ws.Cells.Item(3, 16).Value = 0.02305628730248842
$ws.Cells.Item(3, 17).Value = 297.6441565157503
$ws.Cells.Item(3, 18).Value = 2678.797408641753
$ws.Cells.Item(3, 19).Value = 0.009951420802279797
$ws.Cells.Item(3, 20).Value = 0.01391926492192939

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 88.281043
$ws.Cells.Item(4, 8).Value = 264.843129
$ws.Cells.Item(4, 9).Value = 0.5104595351890647
$ws.Cells.Item(4, 10).Value = 0.6037079925017727
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 0.6617113333333333
$ws.Cells.Item(4, 14).Value = 1.985134
$ws.Cells.Item(4, 15).Value = 0.003826153930268694
$ws.Cells.Item(4, 16).Value = 0.00452509856122042
$ws.Cells.Item(4, 17).Value = 58.41656667158733
$ws.Cells.Item(4, 18).Value = 525.7491000442859
$ws.Cells.Item(4, 19).Value = 0.001953096756806771
$ws.Cells.Item(4, 20).Value = 0.00273183816826704

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 88.281043
$ws.Cells.Item(5, 8).Value = 264.843129
$ws.Cells.Item(5, 9).Value = 0.5104595351890647
$ws.Cells.Item(5, 10).Value = 0.6037079925017727
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 0.4912783333333333
$ws.Cells.Item(5, 14).Value = 1.473835
$ws.Cells.Item(5, 15).Value = 0.00284067452263553
$ws.Cells.Item(5, 16).Value = 0.003359596197524347
$ws.Cells.Item(5, 17).Value = 43.37056366996833
$ws.Cells.Item(5, 18).Value = 390.335073029715
$ws.Cells.Item(5, 19).Value = 0.001450049396447951
$ws.Cells.Item(5, 20).Value = 0.002028215076024012

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 88.281043
$ws.Cells.Item(6, 8).Value = 264.843129
$ws.Cells.Item(6, 9).Value = 0.5104595351890647
$ws.Cells.Item(6, 10).Value = 0.6037079925017727
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 13).Value = 80.138668
$ws.Cells.Item(6, 14).Value = 160.277336
$ws.Cells.Item(6, 15).Value = 0.4633786125289749
$ws.Cells.Item(6, 16).Value = 0.365351025436994
$ws.Cells.Item(6, 17).Value = 7074.725195670723
$ws.Cells.Item(6, 18).Value = 42448.35117402434
$ws.Cells.Item(6, 19).Value = 0.2365360311680942
$ws.Cells.Item(6, 20).Value = 0.2205653341250318

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 3.371552333333333
$ws.Cells.Item(7, 8).Value = 10.114657
$ws.Cells.Item(7, 9).Value = 0.01949502382905625
$ws.Cells.Item(7, 10).Value = 0.02305628730248842
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 88.281043
$ws.Cells.Item(7, 14).Value = 264.843129
$ws.Cells.Item(7, 15).Value = 0.5104595351890647
$ws.Cells.Item(7, 16).Value = 0.6037079925017727
$ws.Cells.Item(7, 17).Value = 297.6441565157503
$ws.Cells.Item(7, 18).Value = 2678.797408641753
$ws.Cells.Item(7, 19).Value = 0.009951420802279797
$ws.Cells.Item(7, 20).Value = 0.01391926492192939

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 3.371552333333333
$ws.Cells.Item(8, 8).Value = 10.114657
$ws.Cells.Item(8, 9).Value = 0.01949502382905625
$ws.Cells.Item(8, 10).Value = 0.02305628730248842
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 3.371552333333333
$ws.Cells.Item(8, 14).Value = 10.114657
$ws.Cells.Item(8, 15).Value = 0.01949502382905625
$ws.Cells.Item(8, 16).Value = 0.02305628730248842
$ws.Cells.Item(8, 17).Value = 11.36736513640544
$ws.Cells.Item(8, 18).Value = 102.306286227649
$ws.Cells.Item(8, 19).Value = 0.0003800559540954712
$ws.Cells.Item(8, 20).Value = 0.0005315923841748886

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 3.371552333333333
$ws.Cells.Item(9, 8).Value = 10.114657
$ws.Cells.Item(9, 9).Value = 0.01949502382905625
$ws.Cells.Item(9, 10).Value = 0.02305628730248842
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 0.6617113333333333
$ws.Cells.Item(9, 14).Value = 1.985134
$ws.Cells.Item(9, 15).Value = 0.003826153930268694
$ws.Cells.Item(9, 16).Value = 0.00452509856122042
$ws.Cells.Item(9, 17).Value = 2.230994389893111
$ws.Cells.Item(9, 18).Value = 20.078949509038
$ws.Cells.Item(9, 19).Value = 0.00007459096204422544
$ws.Cells.Item(9, 20).Value = 0.000104331972499575

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 3.371552333333333
$ws.Cells.Item(10, 8).Value = 10.114657
$ws.Cells.Item(10, 9).Value = 0.01949502382905625
$ws.Cells.Item(10, 10).Value = 0.02305628730248842
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 0.4912783333333333
$ws.Cells.Item(10, 14).Value = 1.473835
$ws.Cells.Item(10, 15).Value = 0.00284067452263553
$ws.Cells.Item(10, 16).Value = 0.003359596197524347
$ws.Cells.Item(10, 17).Value = 1.656370611066111
$ws.Cells.Item(10, 18).Value = 14.907335499595
$ws.Cells.Item(10, 19).Value = 0.00005537901750937267
$ws.Cells.Item(10, 20).Value = 0.00007745981515046898

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 3.371552333333333
$ws.Cells.Item(11, 8).Value = 10.114657
$ws.Cells.Item(11, 9).Value = 0.01949502382905625
$ws.Cells.Item(11, 10).Value = 0.02305628730248842
$ws.Cells.Item(11, 11).Value = 2
$ws.Cells.Item(11, 13).Value = 80.138668
$ws.Cells.Item(11, 14).Value = 160.277336
$ws.Cells.Item(11, 15).Value = 0.4633786125289749
$ws.Cells.Item(11, 16).Value = 0.365351025436994
$ws.Cells.Item(11, 17).Value = 270.1917130856253
$ws.Cells.Item(11, 18).Value = 1621.150278513752
$ws.Cells.Item(11, 19).Value = 0.00903357709312739
$ws.Cells.Item(11, 20).Value = 0.00842363820873409

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 0.6617113333333333
$ws.Cells.Item(12, 8).Value = 1.985134
$ws.Cells.Item(12, 9).Value = 0.003826153930268694
$ws.Cells.Item(12, 10).Value = 0.00452509856122042
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 88.281043
$ws.Cells.Item(12, 14).Value = 264.843129
$ws.Cells.Item(12, 15).Value = 0.5104595351890647
$ws.Cells.Item(12, 16).Value = 0.6037079925017727
$ws.Cells.Item(12, 17).Value = 58.41656667158733
$ws.Cells.Item(12, 18).Value = 525.7491000442859
$ws.Cells.Item(12, 19).Value = 0.001953096756806771
$ws.Cells.Item(12, 20).Value = 0.00273183816826704

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 0.6617113333333333
$ws.Cells.Item(13, 8).Value = 1.985134
$ws.Cells.Item(13, 9).Value = 0.003826153930268694
$ws.Cells.Item(13, 10).Value = 0.00452509856122042
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 3.371552333333333
$ws.Cells.Item(13, 14).Value = 10.114657
$ws.Cells.Item(13, 15).Value = 0.01949502382905625
$ws.Cells.Item(13, 16).Value = 0.02305628730248842
$ws.Cells.Item(13, 17).Value = 2.230994389893111
$ws.Cells.Item(13, 18).Value = 20.078949509038
$ws.Cells.Item(13, 19).Value = 0.00007459096204422544
$ws.Cells.Item(13, 20).Value = 0.000104331972499575

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 0.6617113333333333
$ws.Cells.Item(14, 8).Value = 1.985134
$ws.Cells.Item(14, 9).Value = 0.003826153930268694
$ws.Cells.Item(14, 10).Value = 0.00452509856122042
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 0.6617113333333333
$ws.Cells.Item(14, 14).Value = 1.985134
$ws.Cells.Item(14, 15).Value = 0.003826153930268694
$ws.Cells.Item(14, 16).Value = 0.00452509856122042
$ws.Cells.Item(14, 17).Value = 0.4378618886617778
$ws.Cells.Item(14, 18).Value = 3.940756997956
$ws.Cells.Item(14, 19).Value = 0.00001463945389811058
$ws.Cells.Item(14, 20).Value = 0.00002047651698875911

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 0.6617113333333333
$ws.Cells.Item(15, 8).Value = 1.985134
$ws.Cells.Item(15, 9).Value = 0.003826153930268694
$ws.Cells.Item(15, 10).Value = 0.00452509856122042
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 0.4912783333333333
$ws.Cells.Item(15, 14).Value = 1.473835
$ws.Cells.Item(15, 15).Value = 0.00284067452263553
$ws.Cells.Item(15, 16).Value = 0.003359596197524347
$ws.Cells.Item(15, 17).Value = 0.3250844409877778
$ws.Cells.Item(15, 18).Value = 2.92575996889
$ws.Cells.Item(15, 19).Value = 0.00001086885798939608
$ws.Cells.Item(15, 20).Value = 0.00001520250391969902

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 0.6617113333333333
$ws.Cells.Item(16, 8).Value = 1.985134
$ws.Cells.Item(16, 9).Value = 0.003826153930268694
$ws.Cells.Item(16, 10).Value = 0.00452509856122042
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 13).Value = 80.138668
$ws.Cells.Item(16, 14).Value = 160.277336
$ws.Cells.Item(16, 15).Value = 0.4633786125289749
$ws.Cells.Item(16, 16).Value = 0.365351025436994
$ws.Cells.Item(16, 17).Value = 53.02866485383733
$ws.Cells.Item(16, 18).Value = 318.171989123024
$ws.Cells.Item(16, 19).Value = 0.001772957899530192
$ws.Cells.Item(16, 20).Value = 0.001653249399545347

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 0.4912783333333333
$ws.Cells.Item(17, 8).Value = 1.473835
$ws.Cells.Item(17, 9).Value = 0.00284067452263553
$ws.Cells.Item(17, 10).Value = 0.003359596197524347
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 88.281043
$ws.Cells.Item(17, 14).Value = 264.843129
$ws.Cells.Item(17, 15).Value = 0.5104595351890647
$ws.Cells.Item(17, 16).Value = 0.6037079925017727
$ws.Cells.Item(17, 17).Value = 43.37056366996833
$ws.Cells.Item(17, 18).Value = 390.335073029715
$ws.Cells.Item(17, 19).Value = 0.001450049396447951
$ws.Cells.Item(17, 20).Value = 0.002028215076024012

$ws.Cells.Item(18, 5).Value = 3
$ws.Cells.Item(18, 7).Value = 0.4912783333333333
$ws.Cells.Item(18, 8).Value = 1.473835
$ws.Cells.Item(18, 9).Value = 0.00284067452263553
$ws.Cells.Item(18, 10).Value = 0.003359596197524347
$ws.Cells.Item(18, 11).Value = 3
$ws.Cells.Item(18, 13).Value = 3.371552333333333
$ws.Cells.Item(18, 14).Value = 10.114657
$ws.Cells.Item(18, 15).Value = 0.01949502382905625
$ws.Cells.Item(18, 16).Value = 0.02305628730248842
$ws.Cells.Item(18, 17).Value = 1.656370611066111
$ws.Cells.Item(18, 18).Value = 14.907335499595
$ws.Cells.Item(18, 19).Value = 0.00005537901750937267
$ws.Cells.Item(18, 20).Value = 0.00007745981515046898

$ws.Cells.Item(19, 5).Value = 3
$ws.Cells.Item(19, 7).Value = 0.4912783333333333
$ws.Cells.Item(19, 8).Value = 1.473835
$ws.Cells.Item(19, 9).Value = 0.00284067452263553
$ws.Cells.Item(19, 10).Value = 0.003359596197524347
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 13).Value = 0.6617113333333333
$ws.Cells.Item(19, 14).Value = 1.985134
$ws.Cells.Item(19, 15).Value = 0.003826153930268694
$ws.Cells.Item(19, 16).Value = 0.00452509856122042
$ws.Cells.Item(19, 17).Value = 0.3250844409877778
$ws.Cells.Item(19, 18).Value = 2.92575996889
$ws.Cells.Item(19, 19).Value = 0.00001086885798939608
$ws.Cells.Item(19, 20).Value = 0.00001520250391969902

$ws.Cells.Item(20, 5).Value = 3
$ws.Cells.Item(20, 7).Value = 0.4912783333333333
$ws.Cells.Item(20, 8).Value = 1.473835
$ws.Cells.Item(20, 9).Value = 0.00284067452263553
$ws.Cells.Item(20, 10).Value = 0.003359596197524347
$ws.Cells.Item(20, 11).Value = 3
$ws.Cells.Item(20, 13).Value = 0.4912783333333333
$ws.Cells.Item(20, 14).Value = 1.473835
$ws.Cells.Item(20, 15).Value = 0.00284067452263553
$ws.Cells.Item(20, 16).Value = 0.003359596197524347
$ws.Cells.Item(20, 17).Value = 0.2413544008027778
$ws.Cells.Item(20, 18).Value = 2.172189607225
$ws.Cells.Item(20, 19).Value = 0.000008069431743550599
$ws.Cells.Item(20, 20).Value = 0.00001128688661042005

$ws.Cells.Item(21, 5).Value = 3
$ws.Cells.Item(21, 7).Value = 0.4912783333333333
$ws.Cells.Item(21, 8).Value = 1.473835
$ws.Cells.Item(21, 9).Value = 0.00284067452263553
$ws.Cells.Item(21, 10).Value = 0.003359596197524347
$ws.Cells.Item(21, 11).Value = 2
$ws.Cells.Item(21, 13).Value = 80.138668
$ws.Cells.Item(21, 14).Value = 160.277336
$ws.Cells.Item(21, 15).Value = 0.4633786125289749
$ws.Cells.Item(21, 16).Value = 0.365351025436994
$ws.Cells.Item(21, 17).Value = 39.37039125059333
$ws.Cells.Item(21, 18).Value = 236.22234750356
$ws.Cells.Item(21, 19).Value = 0.00131630781894526
$ws.Cells.Item(21, 20).Value = 0.001227431915819746

$ws.Cells.Item(22, 5).Value = 2
$ws.Cells.Item(22, 7).Value = 80.138668
$ws.Cells.Item(22, 8).Value = 160.277336
$ws.Cells.Item(22, 9).Value = 0.4633786125289749
$ws.Cells.Item(22, 10).Value = 0.365351025436994
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 13).Value = 88.281043
$ws.Cells.Item(22, 14).Value = 264.843129
$ws.Cells.Item(22, 15).Value = 0.5104595351890647
$ws.Cells.Item(22, 16).Value = 0.6037079925017727
$ws.Cells.Item(22, 17).Value = 7074.725195670723
$ws.Cells.Item(22, 18).Value = 42448.35117402434
$ws.Cells.Item(22, 19).Value = 0.2365360311680942
$ws.Cells.Item(22, 20).Value = 0.2205653341250318

$ws.Cells.Item(23, 5).Value = 2
$ws.Cells.Item(23, 7).Value = 80.138668
$ws.Cells.Item(23, 8).Value = 160.277336
$ws.Cells.Item(23, 9).Value = 0.4633786125289749
$ws.Cells.Item(23, 10).Value = 0.365351025436994
$ws.Cells.Item(23, 11).Value = 3
$ws.Cells.Item(23, 13).Value = 3.371552333333333
$ws.Cells.Item(23, 14).Value = 10.114657
$ws.Cells.Item(23, 15).Value = 0.01949502382905625
$ws.Cells.Item(23, 16).Value = 0.02305628730248842
$ws.Cells.Item(23, 17).Value = 270.1917130856253
$ws.Cells.Item(23, 18).Value = 1621.150278513752
$ws.Cells.Item(23, 19).Value = 0.00903357709312739
$ws.Cells.Item(23, 20).Value = 0.00842363820873409

$ws.Cells.Item(24, 5).Value = 2
$ws.Cells.Item(24, 7).Value = 80.138668
$ws.Cells.Item(24, 8).Value = 160.277336
$ws.Cells.Item(24, 9).Value = 0.4633786125289749
$ws.Cells.Item(24, 10).Value = 0.365351025436994
$ws.Cells.Item(24, 11).Value = 3
$ws.Cells.Item(24, 13).Value = 0.6617113333333333
$ws.Cells.Item(24, 14).Value = 1.985134
$ws.Cells.Item(24, 15).Value = 0.003826153930268694
$ws.Cells.Item(24, 16).Value = 0.00452509856122042
$ws.Cells.Item(24, 17).Value = 53.02866485383733
$ws.Cells.Item(24, 18).Value = 318.171989123024
$ws.Cells.Item(24, 19).Value = 0.001772957899530192
$ws.Cells.Item(24, 20).Value = 0.001653249399545347

$ws.Cells.Item(25, 5).Value = 2
$ws.Cells.Item(25, 7).Value = 80.138668
$ws.Cells.Item(25, 8).Value = 160.277336
$ws.Cells.Item(25, 9).Value = 0.4633786125289749
$ws.Cells.Item(25, 10).Value = 0.365351025436994
$ws.Cells.Item(25, 11).Value = 3
$ws.Cells.Item(25, 13).Value = 0.4912783333333333
$ws.Cells.Item(25, 14).Value = 1.473835
$ws.Cells.Item(25, 15).Value = 0.00284067452263553
$ws.Cells.Item(25, 16).Value = 0.003359596197524347
$ws.Cells.Item(25, 17).Value = 39.37039125059333
$ws.Cells.Item(25, 18).Value = 236.22234750356
$ws.Cells.Item(25, 19).Value = 0.00131630781894526
$ws.Cells.Item(25, 20).Value = 0.001227431915819746

$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 7).Value = 80.138668
$ws.Cells.Item(26, 8).Value = 160.277336
$ws.Cells.Item(26, 9).Value = 0.4633786125289749
$ws.Cells.Item(26, 10).Value = 0.365351025436994
$ws.Cells.Item(26, 11).Value = 2
$ws.Cells.Item(26, 13).Value = 80.138668
$ws.Cells.Item(26, 14).Value = 160.277336
$ws.Cells.Item(26, 15).Value = 0.4633786125289749
$ws.Cells.Item(26, 16).Value = 0.365351025436994
$ws.Cells.Item(26, 17).Value = 6422.206108814224
$ws.Cells.Item(26, 18).Value = 25688.82443525689
$ws.Cells.Item(26, 19).Value = 0.2147197385492778
$ws.Cells.Item(26, 20).Value = 0.1334813717878631
